$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.809.16'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '2.538.05'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.45'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.41'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.582'
$ws.Range('D9').Value = '2.534.95'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.47'
$ws.Range('E11').Value = '  -5.57%  '
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.14'
$ws.Range('E14').Value = '  -3.38%  '
$ws.Range('D15').Value = '2.990.94'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = '62.754.20'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').Value = '2.529.89'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '333.49'
$ws.Range('E20').Value = '  -3.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  -1.97%  '
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.95'
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.169'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.28'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.45'
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.22'
$ws.Range('E30').Value = '  +5.70%  '
$ws.Range('D31').Value = '0.0₃0804'
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.99'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  -3.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '397.18'
$ws.Range('E35').Value = '  -7.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.98'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  -3.60%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.72'
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.20'
$ws.Range('E42').Value = '  -3.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '150.27'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E44').Value = '  -2.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.54'
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0236'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.97'
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.29'
$ws.Range('E51').Value = '  +0.20%  '
